$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate stats to reflect the newly-closed trade.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.99   # Current Capital
$wsSummary.Range("B4").Value = -0.01     # Total P&L $
$wsSummary.Range("B5").Value = -0.01     # Total P&L %
$wsSummary.Range("B6").Value = 17        # Total Trades
$wsSummary.Range("B8").Value = 7         # Losing Trades
$wsSummary.Range("B9").Value = 35.29     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4).
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.98999999999999   # Capital
$wsStatus.Range("D4").Value = 17                  # Trades
$wsStatus.Range("E4").Value = -0.01                # P&L $
$wsStatus.Range("F4").Value = -0.01                # P&L %
$wsStatus.Range("G4").Value = 35.29               # Win Rate %

# ---------------------------------------------------------------------------
# Helper: append the new closed trade (Trade #17) as row 18 of a trade log
# sheet. Date/Time columns must stay plain text (matching every other row in
# these logs), so we force a Text number format before assigning the value
# and then restore the "Normal" style afterwards so no stray style index is
# left behind on the cell.
# ---------------------------------------------------------------------------
function Add-Trade17Row {
    param($ws)

    $ws.Range("A18").Value = 17

    $ws.Range("B18").NumberFormat = "@"
    $ws.Range("B18").Value = "2026-02-17"
    $ws.Range("B18").Style = "Normal"

    $ws.Range("C18").NumberFormat = "@"
    $ws.Range("C18").Value = "12:28:43"
    $ws.Range("C18").Style = "Normal"

    $ws.Range("D18").Value = "MarketMaking"
    $ws.Range("E18").Value = "UP"
    $ws.Range("F18").Value = 0.07000000000000001
    $ws.Range("G18").Value = 0.051099
    $ws.Range("H18").Value = "CLOSED"
    $ws.Range("I18").Value = -27.002
    $ws.Range("J18").Value = -0.02
    $ws.Range("K18").Value = 99.98999999999999
    $ws.Range("L18").Value = 0
    $ws.Range("M18").Value = 0
    $ws.Range("N18").Value = 0.6
    $ws.Range("O18").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P18").Value = "early_exit"
    $ws.Range("Q18").Value = 0.11
}

# ---------------------------------------------------------------------------
# Sheet "All Trades": append the new trade row.
# ---------------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-Trade17Row $wsAllTrades

# ---------------------------------------------------------------------------
# Sheet "MarketMaking": append the identical new trade row.
# ---------------------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade17Row $wsMarketMaking
